$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells stay as Text so values like "0.100" or "42.655.71" are not
# reinterpreted as numbers/dates by Excel when assigned via .Value

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.655.71'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.267.66'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.634'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '77.07'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +6.72%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.638'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.97'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.24'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.605.40'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.91'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.861'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.266.43'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.566.61'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0989'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.11'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '234.40'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.77'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -5.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.26'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '167.62'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.88'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.40'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0853'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +6.19%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.77%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '31.02'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.57'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.70'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0304'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '13.75'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +7.69%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.26'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.85'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.207'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '109.16'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +14.20%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '61.15'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -4.27%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.64'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -7.74%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.100'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -3.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0000339'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +131.07%  '
